$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("userinfo")

$data = @(
    @("FirstName", "LastName", "Email", "Phone", "Address", "Company Name", "Role in Company", "status"),
    @("Roger", "Huffman", "RogerLHuffman@armyspy.com ", "619-675-5903", "2299 Grim Avenue. San Diego, CA 92103 ", "Heilig-Meyers", "Floor finisher", "success"),
    @("Vesta", "McAteer", "VestaEMcAteer@jourrapide.com ", "253-520-6265", "115 Horizon Circle. Kent, WA 98031 ", "Disc Jockey", "Telecommunications specialist", "success"),
    @("Laurie", "Strand", "LaurieLStrand@rhyta.com ", "901-870-5509", "3130 Gladwell Street. Memphis, TN 38133 ", "Total Yard Management", "Loss prevention agent", "success")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
